# Update expected validation data to reflect 5% mobilization cost.
# (The develop branch had mobilization cost at 10% of total module cost
# for utility-scale Foundation modules; this changes it to 5%, halving
# the previously-computed FoundationCost/Mobilization cost figures, and
# re-points the sheet's AutoFilter at the FoundationCost/Mobilization
# rows so they're the ones left visible.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Halve the FoundationCost/Mobilization "Cost per turbine" (G) and
# "Cost per project" (H) figures for the four utility-scale validation
# rows that carry nonzero costs.
$ws.Range("G113").Value = 4779.1223863116402
$ws.Range("H113").Value = 477912.23863116401

$ws.Range("G146").Value = 6883.5568677524498
$ws.Range("H146").Value = 413013.412065148

$ws.Range("G179").Value = 11297.25833822005
$ws.Range("H179").Value = 316323.23347016203

$ws.Range("G212").Value = 17047.588862193352
$ws.Range("H212").Value = 477332.48814141349

# Re-point the sheet's AutoFilter: Module (col D) = FoundationCost AND
# Type of cost (col F) = Mobilization, instead of the old Module =
# ErectionCost criterion. (Using the xlFilterValues form so the saved
# XML uses a plain <filters><filter val="..."/></filters> list, matching
# how Excel stores a simple checkbox-list AutoFilter selection.)
$ws.Range("A1:I241").AutoFilter(4, @("FoundationCost"), 7)
$ws.Range("A1:I241").AutoFilter(6, @("Mobilization"), 7)

# Make sure every data row's hidden state matches the new filter
# criteria exactly (belt-and-suspenders around the filter re-apply).
for ($r = 2; $r -le 241; $r++) {
    $moduleVal = $ws.Cells.Item($r, 4).Value2
    $costTypeVal = $ws.Cells.Item($r, 6).Value2
    if ($moduleVal -eq "FoundationCost" -and $costTypeVal -eq "Mobilization") {
        $ws.Rows.Item($r).Hidden = $false
    } else {
        $ws.Rows.Item($r).Hidden = $true
    }
}

# Reflect where the user ended up after filtering/scrolling.
$ws.Range("K246").Select()
